# Update "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column between "Week" (A) and "ASIN" (B)
#  - populate it with the week-start dates (kept as literal text, not real dates)
#  - drop the leading zero in the "Week" labels (W01 -> W1, ... W09 -> W9)
#  - convert the "is_holiday_week" column (now column J) to real booleans

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- insert the new column, shifting ASIN..is_holiday_week one column right ---
$ws.Columns.Item(2).Insert()

# --- header ---
$ws.Cells.Item(1, 2).Value2 = "Week_Start_Date"

# --- week start dates, keyed by row ---
$weekStartDates = @{
    2  = "2024-12-15"
    3  = "2024-12-22"
    4  = "2024-12-29"
    5  = "2025-01-05"
    6  = "2025-01-12"
    7  = "2025-01-19"
    8  = "2025-01-26"
    9  = "2025-02-02"
    10 = "2025-02-09"
    11 = "2025-02-16"
    12 = "2025-02-23"
    13 = "2025-03-02"
    14 = "2025-03-09"
    15 = "2025-03-16"
    16 = "2025-03-23"
    17 = "2025-03-30"
}

for ($r = 2; $r -le 17; $r++) {
    # keep the value as literal text, not an auto-converted date serial
    $cell = $ws.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value2 = $weekStartDates[$r]

    # strip the leading zero from the week label: W01 -> W1 ... W09 -> W9
    $weekCell = $ws.Cells.Item($r, 1)
    $week = $weekCell.Value2
    if ($week -match '^W0(\d)$') {
        $weekCell.Value2 = "W" + $matches[1]
    }

    # is_holiday_week now lives in column J -- make it a real boolean
    $holidayCell = $ws.Cells.Item($r, 10)
    $holidayCell.Value2 = [bool]([int]$holidayCell.Value2)
}
